# This edit reorders/rotates the 31 data rows (rows 2..32) of the sheet.
# Every row keeps its full original content (all 20 columns), but the
# content is moved to a different row position according to the mapping
# below (new row number -> original/source row number). No cell values
# are actually altered - only their row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 32
$lastCol = 20  # column T

# 1) Snapshot all current values for rows 2..32, columns A..T (1..20),
#    before writing anything back, so source data is never clobbered
#    mid-update.
$data = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2())
    }
    $data[$r] = $rowVals
}

# 2) Mapping of destination row -> source row describing where each
#    row's data should end up.
$map = @{
    2 = 4;  3 = 6;  4 = 22; 5 = 25; 6 = 20; 7 = 17; 8 = 9;  9 = 18; 10 = 24;
    11 = 27; 12 = 28; 13 = 32; 14 = 23; 15 = 13; 16 = 14; 17 = 29; 18 = 2;
    19 = 3; 20 = 7; 21 = 8; 22 = 15; 23 = 31; 24 = 21; 25 = 10; 26 = 11;
    27 = 16; 28 = 30; 29 = 19; 30 = 26; 31 = 12; 32 = 5
}

# 3) Write the values back out into their new row positions.
foreach ($newRow in $map.Keys) {
    $srcRow = $map[$newRow]
    $srcVals = $data[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $srcVals[$c - 1]
    }
}
